$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value2 = "'Spettacoli,Teatro"
$ws.Range("B11").Value2 = "'Modena"
$ws.Range("C11").Value2 = "'Viale Caduti in Guerra, 196"
$ws.Range("D11").Value2 = "'2020-09-17T12:45:54+00:00"
$ws.Range("E11").Value2 = "'A cura delle allieve attrici e degli allievi attori di ERT / Teatro Nazionale"
$ws.Range("F11").Value2 = "'2014-09-30T12:50:00+00:00"
$ws.Range("G11").Value2 = "'info@emiliaromagnateatro.com"
$ws.Range("H11").Value2 = "'2022-06-07T22:00:00+00:00"
$ws.Range("I11").Value2 = "'2022-06-11T21:55:00+00:00"
$ws.Range("J11").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/hereafter-episodi-teatrali/@@images/59c6e256-ecc2-48df-912f-b171b831bd83.jpeg"
$ws.Range("K11").Value2 = "'"
$ws.Range("L11").Value2 = "'2022-06-06T07:48:28+00:00"
$ws.Range("M11").Value2 = "'Teatro Tempio"
$ws.Range("N11").Value2 = "' ore 19.00"
$ws.Range("O11").Value2 = "'"
$ws.Range("P11").Value2 = "' A pagamento, vedi nel testo le info per l'acquisto dei biglietti."
$ws.Range("Q11").Value2 = "'"
$ws.Range("R11").Value2 = "'059/2163021"
$ws.Range("S11").Value2 = "'HEREAFTER. Episodi teatrali"
$ws.Range("T11").Value2 = "'"
$ws.Range("U11").Value2 = "'http://www.emiliaromagnateatro.com"
$ws.Range("V11").Value2 = $false
$ws.Range("W11").Value2 = "'"
$ws.Range("X11").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/hereafter-episodi-teatrali"
$ws.Range("Y11").Value2 = "'44,64381951149482"
$ws.Range("Z11").Value2 = "'10,93139345085676"
$ws.Range("AA11").Value2 = "'POINT (10.93139345085676 44.64381951149482)"
$ws.Range("A12").Value2 = "'Conferenze, Seminari, Incontri e Lezioni,Libri"
$ws.Range("B12").Value2 = "'Modena"
$ws.Range("C12").Value2 = "'Corso Vittorio Emanuele, 59"
$ws.Range("D12").Value2 = "'2022-06-04T07:59:44+00:00"
$ws.Range("E12").Value2 = "'Presentazione del libro"
$ws.Range("F12").Value2 = "'2022-06-04T07:59:52+00:00"
$ws.Range("G12").Value2 = "'​info@accademiasla-mo.it"
$ws.Range("H12").Value2 = "'2022-06-08T07:00:00+00:00"
$ws.Range("I12").Value2 = "'2022-06-08T08:00:00+00:00"
$ws.Range("J12").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/delitti-in-prima-pagina-la-giustizia-nella-societa-dell2019informazione/@@images/1a4df534-f8d8-4ac6-87d4-3d9176cd882c.jpeg"
$ws.Range("K12").Value2 = "'"
$ws.Range("L12").Value2 = "'2022-06-04T07:59:52+00:00"
$ws.Range("M12").Value2 = "'Accademia Nazionale di Scienze Lettere e Arti di Modena"
$ws.Range("N12").Value2 = "' ore 15.30"
$ws.Range("O12").Value2 = "'"
$ws.Range("P12").Value2 = "'"
$ws.Range("Q12").Value2 = "'"
$ws.Range("R12").Value2 = "'059 225566"
$ws.Range("S12").Value2 = "'`"Delitti in prima pagina. La giustizia nella società dell’informazione`""
$ws.Range("T12").Value2 = "'"
$ws.Range("U12").Value2 = "'www.accademiasla-mo.it"
$ws.Range("V12").Value2 = $false
$ws.Range("W12").Value2 = 41123
$ws.Range("X12").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/delitti-in-prima-pagina-la-giustizia-nella-societa-dell2019informazione"
$ws.Range("Y12").Value2 = "'44,64582"
$ws.Range("Z12").Value2 = "'10,92572"
$ws.Range("AA12").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A13").Value2 = "'Mostre"
$ws.Range("B13").Value2 = "'Modena"
$ws.Range("C13").Value2 = "'Strada Vaciglio Nord, 6"
$ws.Range("D13").Value2 = "'2022-06-04T08:30:34+00:00"
$ws.Range("E13").Value2 = "'"
$ws.Range("F13").Value2 = "'2022-06-04T08:30:59+00:00"
$ws.Range("G13").Value2 = "'"
$ws.Range("H13").Value2 = "'2022-06-08T08:00:00+00:00"
$ws.Range("I13").Value2 = "'2022-07-08T09:00:00+00:00"
$ws.Range("J13").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/oro-rosso-fragole-pomodori-molestie-e-sfruttamento-nel-mediterraneo/@@images/04c6eef6-5450-4d2e-a11b-c353a7bdb6b0.jpeg"
$ws.Range("K13").Value2 = "'"
$ws.Range("L13").Value2 = "'2022-06-06T10:06:47+00:00"
$ws.Range("M13").Value2 = "'Sala Renata Bergonzoni della Casa delle Donne"
$ws.Range("N13").Value2 = "' Inaugurazione mercoledì 8 giugno ore 18.30  mostra aperta dal 10 giugno nei seguenti orari:  venerdì e sabato dalle 10 alle 13 (ad esclusione di venerdì 17 giugno) "
$ws.Range("O13").Value2 = "'"
$ws.Range("P13").Value2 = "' ingresso libero"
$ws.Range("Q13").Value2 = "'"
$ws.Range("R13").Value2 = "'"
$ws.Range("S13").Value2 = "'Oro rosso. Fragole, pomodori, molestie e sfruttamento nel Mediterraneo"
$ws.Range("T13").Value2 = "'"
$ws.Range("U13").Value2 = "'"
$ws.Range("V13").Value2 = $false
$ws.Range("W13").Value2 = 41123
$ws.Range("X13").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/oro-rosso-fragole-pomodori-molestie-e-sfruttamento-nel-mediterraneo"
$ws.Range("Y13").Value2 = "'44,64582"
$ws.Range("Z13").Value2 = "'10,92572"
$ws.Range("AA13").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A14").Value2 = "'Altri eventi,Musica"
$ws.Range("B14").Value2 = "'Modena"
$ws.Range("C14").Value2 = "'Centro storico"
$ws.Range("D14").Value2 = "'2022-05-20T10:02:04+00:00"
$ws.Range("E14").Value2 = "'Rassegna musicale nella suggestiva piazzetta della Pomposa"
$ws.Range("F14").Value2 = "'2022-05-20T10:02:34+00:00"
$ws.Range("G14").Value2 = "'"
$ws.Range("H14").Value2 = "'2022-05-12T09:00:00+00:00"
$ws.Range("I14").Value2 = "'2022-06-30T10:00:00+00:00"
$ws.Range("J14").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/serate-in-pomposa-1deg-edizione/@@images/90c7f1e5-b97b-4472-866b-68e21ff4b305.jpeg"
$ws.Range("K14").Value2 = "'"
$ws.Range("L14").Value2 = "'2022-05-20T10:02:34+00:00"
$ws.Range("M14").Value2 = "'Piazzetta Pomposa"
$ws.Range("N14").Value2 = "' tutti i mercoledì e giovedì, dalle ore 18.30"
$ws.Range("O14").Value2 = "'"
$ws.Range("P14").Value2 = "'"
$ws.Range("Q14").Value2 = "'"
$ws.Range("R14").Value2 = "'"
$ws.Range("S14").Value2 = "'Serate in Pomposa - 1° edizione"
$ws.Range("T14").Value2 = "'"
$ws.Range("U14").Value2 = "'"
$ws.Range("V14").Value2 = $false
$ws.Range("W14").Value2 = 41123
$ws.Range("X14").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/serate-in-pomposa-1deg-edizione"
$ws.Range("Y14").Value2 = "'44,64582"
$ws.Range("Z14").Value2 = "'10,92572"
$ws.Range("AA14").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A15").Value2 = "'Mostre"
$ws.Range("B15").Value2 = "'Modena"
$ws.Range("C15").Value2 = "'Corso Cavour, angolo corso Canalgrnde"
$ws.Range("D15").Value2 = "'2021-09-09T10:12:02+00:00"
$ws.Range("E15").Value2 = "'"
$ws.Range("F15").Value2 = "'2021-09-09T10:12:26+00:00"
$ws.Range("G15").Value2 = "'info@fmav.org - biglietteria@fmav.org"
$ws.Range("H15").Value2 = "'2022-06-08T10:00:00+00:00"
$ws.Range("I15").Value2 = "'2022-09-18T11:00:00+00:00"
$ws.Range("J15").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/candice-breitz-never-ending-stories/@@images/d5ebf4a4-66cb-44cb-9e8c-79f9fcfc72ea.jpeg"
$ws.Range("K15").Value2 = "'"
$ws.Range("L15").Value2 = "'2022-05-21T09:18:35+00:00"
$ws.Range("M15").Value2 = "'FMAV - Palazzina dei Giardini"
$ws.Range("N15").Value2 = "' dal mercoledì alla domenica dalle 15 alle 19 "
$ws.Range("O15").Value2 = "'"
$ws.Range("P15").Value2 = "' Ingresso 6 € / riduzioni 4 € (Circuito Vivaticket) Ingresso libero: ogni mercoledì "
$ws.Range("Q15").Value2 = "'"
$ws.Range("R15").Value2 = "'059 2033166"
$ws.Range("S15").Value2 = "'Candice Breitz: Never Ending Stories"
$ws.Range("T15").Value2 = "'"
$ws.Range("U15").Value2 = "'https://www.fmav.org"
$ws.Range("V15").Value2 = $false
$ws.Range("W15").Value2 = "'41121"
$ws.Range("X15").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/candice-breitz-never-ending-stories"
$ws.Range("Y15").Value2 = "'44,64582"
$ws.Range("Z15").Value2 = "'10,92572"
$ws.Range("AA15").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A16").Value2 = "'Altri eventi,Iniziative per bambini"
$ws.Range("B16").Value2 = "'Modena"
$ws.Range("C16").Value2 = "'via Selmi, 63"
$ws.Range("D16").Value2 = "'2022-05-17T07:24:22+00:00"
$ws.Range("E16").Value2 = "'"
$ws.Range("F16").Value2 = "'2022-05-17T07:25:56+00:00"
$ws.Range("G16").Value2 = "'"
$ws.Range("H16").Value2 = "'2022-06-10T07:00:00+00:00"
$ws.Range("I16").Value2 = "'2022-06-10T08:00:00+00:00"
$ws.Range("J16").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/laboratori-artistici-in-occasione-della-mostra-dedicata-a-gianni-valbonesi/@@images/c62a264a-bdfb-469d-85eb-8d034a4afbbc.jpeg"
$ws.Range("K16").Value2 = "'"
$ws.Range("L16").Value2 = "'2022-05-26T07:16:41+00:00"
$ws.Range("M16").Value2 = "'Complesso San Paolo"
$ws.Range("N16").Value2 = "' dalle ore 17 alle 18.30"
$ws.Range("O16").Value2 = "'"
$ws.Range("P16").Value2 = "' Iscrizione su prenotazione, è richiesto un contributo di 3€ a partecipante.  "
$ws.Range("Q16").Value2 = "'"
$ws.Range("R16").Value2 = "'"
$ws.Range("S16").Value2 = "'Laboratori artistici in occasione della mostra dedicata a Gianni Valbonesi"
$ws.Range("T16").Value2 = "'"
$ws.Range("U16").Value2 = "'"
$ws.Range("V16").Value2 = $false
$ws.Range("W16").Value2 = 41123
$ws.Range("X16").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/laboratori-artistici-in-occasione-della-mostra-dedicata-a-gianni-valbonesi"
$ws.Range("Y16").Value2 = "'44,64582"
$ws.Range("Z16").Value2 = "'10,92572"
$ws.Range("AA16").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A17").Value2 = "'Conferenze, Seminari, Incontri e Lezioni"
$ws.Range("B17").Value2 = "'Modena"
$ws.Range("C17").Value2 = "'Strada Vaciglio Nord, 6"
$ws.Range("D17").Value2 = "'2022-06-04T08:30:34+00:00"
$ws.Range("E17").Value2 = "'"
$ws.Range("F17").Value2 = "'2022-06-04T08:30:59+00:00"
$ws.Range("G17").Value2 = "'"
$ws.Range("H17").Value2 = "'2022-06-10T08:00:00+00:00"
$ws.Range("I17").Value2 = "'2022-06-10T09:00:00+00:00"
$ws.Range("J17").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/da-suddite-a-cittadine-gabriella-degli-esposti-e-le-partigiane-modenesi-nelle-fonti-documentarie-e-nelle-memorie/@@images/9a9f2a19-02f6-4cd6-a272-5cac5ae7adbd.jpeg"
$ws.Range("K17").Value2 = "'"
$ws.Range("L17").Value2 = "'2022-06-04T08:52:51+00:00"
$ws.Range("M17").Value2 = "'Sala Renata Bergonzoni della Casa delle Donne"
$ws.Range("N17").Value2 = "' ore 18.30"
$ws.Range("O17").Value2 = "'"
$ws.Range("P17").Value2 = "' ingresso libero"
$ws.Range("Q17").Value2 = "'"
$ws.Range("R17").Value2 = "'"
$ws.Range("S17").Value2 = "'Da suddite a cittadine. Gabriella Degli Esposti e le partigiane modenesi nelle fonti documentarie e nelle memorie"
$ws.Range("T17").Value2 = "'"
$ws.Range("U17").Value2 = "'"
$ws.Range("V17").Value2 = $false
$ws.Range("W17").Value2 = 41123
$ws.Range("X17").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/da-suddite-a-cittadine-gabriella-degli-esposti-e-le-partigiane-modenesi-nelle-fonti-documentarie-e-nelle-memorie"
$ws.Range("Y17").Value2 = "'44,64582"
$ws.Range("Z17").Value2 = "'10,92572"
$ws.Range("AA17").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A18").Value2 = "'Visite guidate"
$ws.Range("B18").Value2 = "'Modena"
$ws.Range("C18").Value2 = "'centro storico"
$ws.Range("D18").Value2 = "'2022-04-29T10:05:46+00:00"
$ws.Range("E18").Value2 = "'Visite straordinarie, dalle 19 alle 22, il venerdì, sabato e domenica, dal 1 maggio al 17 luglio."
$ws.Range("F18").Value2 = "'2022-04-29T10:06:04+00:00"
$ws.Range("G18").Value2 = "'"
$ws.Range("H18").Value2 = "'2022-05-01T09:00:00+00:00"
$ws.Range("I18").Value2 = "'2022-07-17T10:00:00+00:00"
$ws.Range("J18").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/torre-ghirlandina-aperture-serali/@@images/bfdf71f9-ed36-4565-8549-0084dd664317.jpeg"
$ws.Range("K18").Value2 = "'Torre Ghirlandina"
$ws.Range("L18").Value2 = "'2022-04-29T10:06:04+00:00"
$ws.Range("M18").Value2 = "'Piazza Torre"
$ws.Range("N18").Value2 = "' Dal 1 maggio al 17 luglio: tutti i venerdì, sabato e domenica dalle 19 alle 22 "
$ws.Range("O18").Value2 = "'"
$ws.Range("P18").Value2 = "'"
$ws.Range("Q18").Value2 = "'"
$ws.Range("R18").Value2 = "'"
$ws.Range("S18").Value2 = "'Torre Ghirlandina, aperture serali"
$ws.Range("T18").Value2 = "'"
$ws.Range("U18").Value2 = "'"
$ws.Range("V18").Value2 = $false
$ws.Range("W18").Value2 = 41123
$ws.Range("X18").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/torre-ghirlandina-aperture-serali"
$ws.Range("Y18").Value2 = "'44,64582"
$ws.Range("Z18").Value2 = "'10,92572"
$ws.Range("AA18").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A19").Value2 = "'Visite guidate,Altri eventi,Spettacoli"
$ws.Range("B19").Value2 = "'Modena"
$ws.Range("C19").Value2 = "'"
$ws.Range("D19").Value2 = "'2022-05-30T09:56:40+00:00"
$ws.Range("E19").Value2 = "'Una serata speciale incentrata intorno e dentro al Duomo di Modena"
$ws.Range("F19").Value2 = "'2022-05-30T09:57:07+00:00"
$ws.Range("G19").Value2 = "'"
$ws.Range("H19").Value2 = "'2022-06-10T09:00:00+00:00"
$ws.Range("I19").Value2 = "'2022-06-10T10:00:00+00:00"
$ws.Range("J19").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/creature-e-custodi-la-lunga-notte-delle-chiese-7degedizione/@@images/5d98cc67-ebaa-4fc2-a6f2-fbd81bef52b2.jpeg"
$ws.Range("K19").Value2 = "'"
$ws.Range("L19").Value2 = "'2022-05-30T09:57:07+00:00"
$ws.Range("M19").Value2 = "'Duomo di Modena"
$ws.Range("N19").Value2 = "' dalle ore 20.30"
$ws.Range("O19").Value2 = "'"
$ws.Range("P19").Value2 = "' Gratuito"
$ws.Range("Q19").Value2 = "'"
$ws.Range("R19").Value2 = "'"
$ws.Range("S19").Value2 = "'Creature e custodi. La lunga notte delle Chiese - 7°edizione"
$ws.Range("T19").Value2 = "'"
$ws.Range("U19").Value2 = "'"
$ws.Range("V19").Value2 = $false
$ws.Range("W19").Value2 = 41123
$ws.Range("X19").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/creature-e-custodi-la-lunga-notte-delle-chiese-7degedizione"
$ws.Range("Y19").Value2 = "'44,64582"
$ws.Range("Z19").Value2 = "'10,92572"
$ws.Range("AA19").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A20").Value2 = "'Concerti,Spettacoli,Altri eventi"
$ws.Range("B20").Value2 = "'Modena"
$ws.Range("C20").Value2 = "'corso Canalgrande, 85"
$ws.Range("D20").Value2 = "'2022-05-30T11:40:24+00:00"
$ws.Range("E20").Value2 = "'Spettacolo di beneficenza per la Croce rossa internazionale, a favore dei profughi e dei feriti dell’Ucraina"
$ws.Range("F20").Value2 = "'2022-05-30T11:41:23+00:00"
$ws.Range("G20").Value2 = "'"
$ws.Range("H20").Value2 = "'2022-06-10T11:00:00+00:00"
$ws.Range("I20").Value2 = "'2022-06-10T12:00:00+00:00"
$ws.Range("J20").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/dolce-epoque-classical-crossover-concert-in-stile-1880-1925/@@images/cbb45c4a-6948-4f24-a365-a1c1a9b8cd7c.jpeg"
$ws.Range("K20").Value2 = "'"
$ws.Range("L20").Value2 = "'2022-05-30T11:42:41+00:00"
$ws.Range("M20").Value2 = "'Teatro comunale Pavarotti-Freni"
$ws.Range("N20").Value2 = "' ore 21.00"
$ws.Range("O20").Value2 = "'"
$ws.Range("P20").Value2 = "'"
$ws.Range("Q20").Value2 = "'"
$ws.Range("R20").Value2 = "'"
$ws.Range("S20").Value2 = "'Dolce époque. Classical Crossover Concert in stile 1880-1925"
$ws.Range("T20").Value2 = "'"
$ws.Range("U20").Value2 = "'"
$ws.Range("V20").Value2 = $false
$ws.Range("W20").Value2 = 41123
$ws.Range("X20").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/dolce-epoque-classical-crossover-concert-in-stile-1880-1925"
$ws.Range("Y20").Value2 = "'44,64582"
$ws.Range("Z20").Value2 = "'10,92572"
$ws.Range("AA20").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A21").Value2 = "'Concerti"
$ws.Range("B21").Value2 = "'Modena"
$ws.Range("C21").Value2 = "'Strada Pomposiana 292"
$ws.Range("D21").Value2 = "'2022-05-30T09:03:14+00:00"
$ws.Range("E21").Value2 = "'Serate musicali nel polo ambientale di Marzaglia "
$ws.Range("F21").Value2 = "'2022-05-30T09:33:45+00:00"
$ws.Range("G21").Value2 = "'"
$ws.Range("H21").Value2 = "'2022-06-10T19:00:00+00:00"
$ws.Range("I21").Value2 = "'2022-06-10T21:00:00+00:00"
$ws.Range("J21").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/marco-ferri-quartet/@@images/85be16b4-09b5-4d77-958a-0d5460201cf0.jpeg"
$ws.Range("K21").Value2 = "'"
$ws.Range("L21").Value2 = "'2022-05-30T09:33:45+00:00"
$ws.Range("M21").Value2 = "'Fattoria Centofiori"
$ws.Range("N21").Value2 = "' Ore 21.00"
$ws.Range("O21").Value2 = "'"
$ws.Range("P21").Value2 = "' Per informazione sui costi contattare tramite  Whatsapp il numero 3293357131 "
$ws.Range("Q21").Value2 = "'"
$ws.Range("R21").Value2 = "'Whatsapp 3293357131"
$ws.Range("S21").Value2 = "'Marco Ferri Quartet"
$ws.Range("T21").Value2 = "'"
$ws.Range("U21").Value2 = "'www.fattoriacentofiori.it"
$ws.Range("V21").Value2 = $false
$ws.Range("W21").Value2 = "'41123"
$ws.Range("X21").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/marco-ferri-quartet"
$ws.Range("Y21").Value2 = "'44,63693000592147"
$ws.Range("Z21").Value2 = "'10,81076003183179"
$ws.Range("AA21").Value2 = "'POINT (10.81076003183179 44.63693000592147)"
$ws.Range("A22").Value2 = "'Altri eventi,Iniziative per bambini"
$ws.Range("B22").Value2 = "'Modena"
$ws.Range("C22").Value2 = "'via Selmi, 63"
$ws.Range("D22").Value2 = "'2022-05-17T07:24:22+00:00"
$ws.Range("E22").Value2 = "'"
$ws.Range("F22").Value2 = "'2022-05-17T07:25:56+00:00"
$ws.Range("G22").Value2 = "'"
$ws.Range("H22").Value2 = "'2022-05-28T07:00:00+00:00"
$ws.Range("I22").Value2 = "'2022-06-11T08:00:00+00:00"
$ws.Range("J22").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/copy_of_laboratori-artistici-in-occasione-della-mostra-dedicata-a-gianni-valbonesi/@@images/faf8f0b5-e078-4310-9f19-8f71719228e5.jpeg"
$ws.Range("K22").Value2 = "'"
$ws.Range("L22").Value2 = "'2022-05-26T07:20:15+00:00"
$ws.Range("M22").Value2 = "'Complesso San Paolo"
$ws.Range("N22").Value2 = "' sabato 28 maggio e sabato 11 giugno dalle ore 10 alle 12.30"
$ws.Range("O22").Value2 = "'"
$ws.Range("P22").Value2 = "' Iscrizione su prenotazione, è richiesto un contributo di 3€ a partecipante.  "
$ws.Range("Q22").Value2 = "'"
$ws.Range("R22").Value2 = "'"
$ws.Range("S22").Value2 = "'Laboratori artistici in occasione della mostra dedicata a Gianni Valbonesi"
$ws.Range("T22").Value2 = "'"
$ws.Range("U22").Value2 = "'"
$ws.Range("V22").Value2 = $false
$ws.Range("W22").Value2 = 41123
$ws.Range("X22").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/copy_of_laboratori-artistici-in-occasione-della-mostra-dedicata-a-gianni-valbonesi"
$ws.Range("Y22").Value2 = "'44,64582"
$ws.Range("Z22").Value2 = "'10,92572"
$ws.Range("AA22").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A23").Value2 = "'Mostre"
$ws.Range("B23").Value2 = "'Modena"
$ws.Range("C23").Value2 = "'Largo Porta Sant’Agostino, 228"
$ws.Range("D23").Value2 = "'2022-06-04T09:45:49+00:00"
$ws.Range("E23").Value2 = "'mostra fotografica di Francesco Jodice"
$ws.Range("F23").Value2 = "'2022-06-04T09:46:06+00:00"
$ws.Range("G23").Value2 = "'info@agomodena.it"
$ws.Range("H23").Value2 = "'2022-06-11T09:00:00+00:00"
$ws.Range("I23").Value2 = "'2022-08-28T10:00:00+00:00"
$ws.Range("J23").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/ritratti-di-classe/@@images/3dba106b-6f95-4190-991d-b13abf85501a.jpeg"
$ws.Range("K23").Value2 = "'Ritratti di classe"
$ws.Range("L23").Value2 = "'2022-06-04T09:47:46+00:00"
$ws.Range("M23").Value2 = "'AGO Modena Fabbriche culturali"
$ws.Range("N23").Value2 = "' vedi sul sito dell'evento"
$ws.Range("O23").Value2 = "'"
$ws.Range("P23").Value2 = "'"
$ws.Range("Q23").Value2 = "'"
$ws.Range("R23").Value2 = "'"
$ws.Range("S23").Value2 = "'Ritratti di Classe"
$ws.Range("T23").Value2 = "'"
$ws.Range("U23").Value2 = "'www.agomodena.it"
$ws.Range("V23").Value2 = $false
$ws.Range("W23").Value2 = 41123
$ws.Range("X23").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/ritratti-di-classe"
$ws.Range("Y23").Value2 = "'44,64582"
$ws.Range("Z23").Value2 = "'10,92572"
$ws.Range("AA23").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A24").Value2 = "'Altri eventi,Fiere, Esposizioni e Mostre Mercato"
$ws.Range("B24").Value2 = "'Modena"
$ws.Range("C24").Value2 = "'Centro storico"
$ws.Range("D24").Value2 = "'2022-05-21T10:32:36+00:00"
$ws.Range("E24").Value2 = "'"
$ws.Range("F24").Value2 = "'2022-05-21T10:32:42+00:00"
$ws.Range("G24").Value2 = "'"
$ws.Range("H24").Value2 = "'2022-06-11T10:00:00+00:00"
$ws.Range("I24").Value2 = "'2022-06-12T11:00:00+00:00"
$ws.Range("J24").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/nocinopoli-la-citta-del-nocino/@@images/0c419aa8-b36d-411f-ae54-b60ae09f6bd7.jpeg"
$ws.Range("K24").Value2 = "'"
$ws.Range("L24").Value2 = "'2022-05-21T10:36:00+00:00"
$ws.Range("M24").Value2 = "'Piazza Matteotti"
$ws.Range("N24").Value2 = "' sabato dalle ore 11 alle 21  domenica dalle ore 10.30 alle 21"
$ws.Range("O24").Value2 = "'"
$ws.Range("P24").Value2 = "'"
$ws.Range("Q24").Value2 = "'"
$ws.Range("R24").Value2 = "'"
$ws.Range("S24").Value2 = "'Nocinopoli - La città del nocino"
$ws.Range("T24").Value2 = "'"
$ws.Range("U24").Value2 = "'"
$ws.Range("V24").Value2 = $false
$ws.Range("W24").Value2 = 41123
$ws.Range("X24").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/nocinopoli-la-citta-del-nocino"
$ws.Range("Y24").Value2 = "'44,64582"
$ws.Range("Z24").Value2 = "'10,92572"
$ws.Range("AA24").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A25").Value2 = "'Mercati"
$ws.Range("B25").Value2 = "'Modena"
$ws.Range("C25").Value2 = "'parco Novi Sad (ex - Piazza d'Armi)"
$ws.Range("D25").Value2 = "'2019-12-10T09:34:45+00:00"
$ws.Range("E25").Value2 = "'Mercato di ambulanti che propongono prodotti sia alimentari che extralimentari"
$ws.Range("F25").Value2 = "'2021-10-09T08:00:00+00:00"
$ws.Range("G25").Value2 = "'"
$ws.Range("H25").Value2 = "'2022-06-11T22:00:00+00:00"
$ws.Range("I25").Value2 = "'2022-06-12T21:55:00+00:00"
$ws.Range("J25").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/copy4_of_fatto-in-italia/@@images/b8f32f6a-f0bd-4736-bb21-b49871c73eec.jpeg"
$ws.Range("K25").Value2 = "'"
$ws.Range("L25").Value2 = "'2022-04-07T09:39:39+00:00"
$ws.Range("M25").Value2 = "'"
$ws.Range("N25").Value2 = "' dalle 7.00 alle 14.30"
$ws.Range("O25").Value2 = "'"
$ws.Range("P25").Value2 = "'"
$ws.Range("Q25").Value2 = "'"
$ws.Range("R25").Value2 = "'"
$ws.Range("S25").Value2 = "'Fatto in Italia"
$ws.Range("T25").Value2 = "'"
$ws.Range("U25").Value2 = "'http://www.consorzioilmercato.it/"
$ws.Range("V25").Value2 = $false
$ws.Range("W25").Value2 = "'"
$ws.Range("X25").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/copy4_of_fatto-in-italia"
$ws.Range("Y25").Value2 = "'44,65137034620577"
$ws.Range("Z25").Value2 = "'10,921029194828652"
$ws.Range("AA25").Value2 = "'POINT (10.921029194828652 44.65137034620577)"
$ws.Range("A26").Value2 = "'Conferenze, Seminari, Incontri e Lezioni,Libri"
$ws.Range("B26").Value2 = "'Modena"
$ws.Range("C26").Value2 = "'via S.Pietro, 1"
$ws.Range("D26").Value2 = "'2022-05-30T10:02:38+00:00"
$ws.Range("E26").Value2 = "'Presentazione del libro"
$ws.Range("F26").Value2 = "'2022-05-30T10:02:53+00:00"
$ws.Range("G26").Value2 = "'"
$ws.Range("H26").Value2 = "'2022-06-12T09:00:00+00:00"
$ws.Range("I26").Value2 = "'2022-06-12T10:00:00+00:00"
$ws.Range("J26").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/contro-la-guerra/@@images/9575eb12-02b7-4af0-8868-9cf74809a77f.jpeg"
$ws.Range("K26").Value2 = "'"
$ws.Range("L26").Value2 = "'2022-05-30T10:02:53+00:00"
$ws.Range("M26").Value2 = "'Chiostro di San Pietro"
$ws.Range("N26").Value2 = "' ore 15.00"
$ws.Range("O26").Value2 = "'"
$ws.Range("P26").Value2 = "'"
$ws.Range("Q26").Value2 = "'"
$ws.Range("R26").Value2 = "'"
$ws.Range("S26").Value2 = "'Contro la guerra"
$ws.Range("T26").Value2 = "'"
$ws.Range("U26").Value2 = "'"
$ws.Range("V26").Value2 = $false
$ws.Range("W26").Value2 = 41123
$ws.Range("X26").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/contro-la-guerra"
$ws.Range("Y26").Value2 = "'44,64582"
$ws.Range("Z26").Value2 = "'10,92572"
$ws.Range("AA26").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A27").Value2 = "'Mercati,Iniziative per bambini"
$ws.Range("B27").Value2 = "'Modena"
$ws.Range("C27").Value2 = "'via Don Pasquino Fiorenzi, 134"
$ws.Range("D27").Value2 = "'2021-08-06T11:34:26+00:00"
$ws.Range("E27").Value2 = "'"
$ws.Range("F27").Value2 = "'2021-08-06T11:36:33+00:00"
$ws.Range("G27").Value2 = "'"
$ws.Range("H27").Value2 = "'2022-06-12T11:00:00+00:00"
$ws.Range("I27").Value2 = "'2022-06-12T12:00:00+00:00"
$ws.Range("J27").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/copy4_of_mercatino-della-polisportiva-madonnina/@@images/7fbd75c5-a620-4574-9c40-f012be16941a.jpeg"
$ws.Range("K27").Value2 = "'logo del mercatino"
$ws.Range("L27").Value2 = "'2022-02-08T10:09:23+00:00"
$ws.Range("M27").Value2 = "'Polisportiva Madonnina"
$ws.Range("N27").Value2 = "' dalle ore 9 alle 13"
$ws.Range("O27").Value2 = "'"
$ws.Range("P27").Value2 = "'"
$ws.Range("Q27").Value2 = "'"
$ws.Range("R27").Value2 = "'340 2607164 Luisa"
$ws.Range("S27").Value2 = "'Mercatino della Polisportiva Madonnina"
$ws.Range("T27").Value2 = "'"
$ws.Range("U27").Value2 = "'https://www.facebook.com/mercatinopolisportivamadonnina/"
$ws.Range("V27").Value2 = $false
$ws.Range("W27").Value2 = 41123
$ws.Range("X27").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/copy4_of_mercatino-della-polisportiva-madonnina"
$ws.Range("Y27").Value2 = "'44,64582"
$ws.Range("Z27").Value2 = "'10,92572"
$ws.Range("AA27").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A28").Value2 = "'Altri eventi,Iniziative per bambini,Visite guidate"
$ws.Range("B28").Value2 = "'Montale Rangone"
$ws.Range("C28").Value2 = "'Via Vandelli (Statale 12 – Nuova Estense)"
$ws.Range("D28").Value2 = "'2021-04-26T15:03:43+00:00"
$ws.Range("E28").Value2 = "'Riapertura al pubblico tutte le domeniche e festivi dal 3 aprile al 19 giugno 2022"
$ws.Range("F28").Value2 = "'2021-04-26T15:06:06+00:00"
$ws.Range("G28").Value2 = "'museo@parcomontale.it"
$ws.Range("H28").Value2 = "'2022-04-03T14:00:00+00:00"
$ws.Range("I28").Value2 = "'2022-06-19T15:00:00+00:00"
$ws.Range("J28").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/parco-archeologico-della-terramara-di-montale/@@images/f7d0a110-c97e-4787-b356-c6a2116576e3.jpeg"
$ws.Range("K28").Value2 = "'Parco archeologico della Terramara di Montale"
$ws.Range("L28").Value2 = "'2022-03-28T14:43:35+00:00"
$ws.Range("M28").Value2 = "'Parco archeologico della Terramara di Montale"
$ws.Range("N28").Value2 = "' Il Parco è aperto dal 3 aprile al 19 giugno, tutte le domeniche e nei giorni festivi di lunedì 18 aprile, lunedì 25 aprile e giovedì 2 giugno"
$ws.Range("O28").Value2 = "'"
$ws.Range("P28").Value2 = "' 7 euro intero, ridotto 5 euro dai 6 ai 13 anni,  gratuito fino a 5 anni e dai 65 anni  riduzione del 50% ai possessori della fidelity card del Parco"
$ws.Range("Q28").Value2 = "'"
$ws.Range("R28").Value2 = "'Per informazioni e prenotazioni: tel. 335 8136948 dalle 9 alle 13 oppure 059 532020 negli orari di apertura del Parco"
$ws.Range("S28").Value2 = "'Parco archeologico della Terramara di Montale"
$ws.Range("T28").Value2 = "'"
$ws.Range("U28").Value2 = "'http://www.parcomontale.it/it"
$ws.Range("V28").Value2 = $false
$ws.Range("W28").Value2 = "'41050"
$ws.Range("X28").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/parco-archeologico-della-terramara-di-montale"
$ws.Range("Y28").Value2 = "'44,64582"
$ws.Range("Z28").Value2 = "'10,92572"
$ws.Range("AA28").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A29").Value2 = "'Altri eventi,Mercati"
$ws.Range("B29").Value2 = "'Modena"
$ws.Range("C29").Value2 = "'Ingresso da viale Berengario"
$ws.Range("D29").Value2 = "'2020-12-29T15:02:32+00:00"
$ws.Range("E29").Value2 = "'"
$ws.Range("F29").Value2 = "'2020-12-29T15:05:00+00:00"
$ws.Range("G29").Value2 = "'"
$ws.Range("H29").Value2 = "'2022-01-03T14:00:00+00:00"
$ws.Range("I29").Value2 = "'2022-12-26T14:59:00+00:00"
$ws.Range("J29").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/mercato-settimanale-del-lunedi-1/@@images/44cdb524-15a5-405b-89af-f9e36f3adeb6.jpeg"
$ws.Range("K29").Value2 = "'"
$ws.Range("L29").Value2 = "'2021-12-28T10:39:07+00:00"
$ws.Range("M29").Value2 = "'Parco Novi Sad"
$ws.Range("N29").Value2 = "' il lunedì dalle ore 8 alle 14"
$ws.Range("O29").Value2 = "'"
$ws.Range("P29").Value2 = "'"
$ws.Range("Q29").Value2 = "'"
$ws.Range("R29").Value2 = "'"
$ws.Range("S29").Value2 = "'Mercato settimanale del lunedì"
$ws.Range("T29").Value2 = "'"
$ws.Range("U29").Value2 = "'http://www.consorzioilmercato.com/"
$ws.Range("V29").Value2 = $false
$ws.Range("W29").Value2 = "'41121"
$ws.Range("X29").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/mercato-settimanale-del-lunedi-1"
$ws.Range("Y29").Value2 = "'44,64582"
$ws.Range("Z29").Value2 = "'10,92572"
$ws.Range("AA29").Value2 = "'POINT (10.92572 44.64582)"
$ws.Range("A30").Value2 = "'Altri eventi,Spettacoli,Musica,Libri"
$ws.Range("B30").Value2 = "'Modena"
$ws.Range("C30").Value2 = "'Centro storico"
$ws.Range("D30").Value2 = "'2022-06-03T15:42:40+00:00"
$ws.Range("E30").Value2 = "'Rassegna di eventi culturali in Piazza XX Settembre"
$ws.Range("F30").Value2 = "'2022-06-03T15:43:48+00:00"
$ws.Range("G30").Value2 = "'"
$ws.Range("H30").Value2 = "'2022-06-07T15:00:00+00:00"
$ws.Range("I30").Value2 = "'2022-06-28T16:00:00+00:00"
$ws.Range("J30").Value2 = "'https://www.comune.modena.it/api/novita/eventi/2022/note-di-stelle/@@images/78120e75-8036-4a4c-b46e-10cc1ff7fe4c.jpeg"
$ws.Range("K30").Value2 = "'"
$ws.Range("L30").Value2 = "'2022-06-04T06:35:10+00:00"
$ws.Range("M30").Value2 = "'Piazza XX Settembre"
$ws.Range("N30").Value2 = "' 7, 14, 21 e 28 giugno alle ore 19 presentazione di libri  alle ore 21 spettacoli"
$ws.Range("O30").Value2 = "'"
$ws.Range("P30").Value2 = "'"
$ws.Range("Q30").Value2 = "'"
$ws.Range("R30").Value2 = "'"
$ws.Range("S30").Value2 = "'Note di Stelle"
$ws.Range("T30").Value2 = "'"
$ws.Range("U30").Value2 = "'"
$ws.Range("V30").Value2 = $false
$ws.Range("W30").Value2 = 41123
$ws.Range("X30").Value2 = "'https://www.comune.modena.it/novita/eventi/2022/note-di-stelle"
$ws.Range("Y30").Value2 = "'44,64582"
$ws.Range("Z30").Value2 = "'10,92572"
$ws.Range("AA30").Value2 = "'POINT (10.92572 44.64582)"
